# Add a new "503A" indicator column (column L) to the facility_info sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("L1").Value = "503A"

# Rows (2-64) that should be flagged "Y" for 503A
$yRows = @(2, 4, 6, 16, 44, 53, 63)

# Default every data row to "NA", then flip the flagged rows to "Y"
for ($r = 2; $r -le 64; $r++) {
    $ws.Cells.Item($r, 12).Value = "NA"
}
foreach ($r in $yRows) {
    $ws.Cells.Item($r, 12).Value = "Y"
}

# Extend the sheet's AutoFilter range to include the new column
$ws.AutoFilterMode = $false
$ws.Range("A1:L64").AutoFilter() | Out-Null

# Keep the hidden _FilterDatabase defined name in sync with the new range
$filterName = $wb.Names.Item("facility_info!_FilterDatabase")
$filterName.RefersTo = "=facility_info!`$A`$1:`$L`$64"

# Leave the selection where the author left it after filling in the column
$ws.Range("L65").Select() | Out-Null
